$d = $word.ActiveDocument

# Paragraph 2: the "m:commentblock" field paragraph.
# The field-code representation (fldChar begin / instrText runs / fldChar end)
# is replaced by literal text runs spelling out the field braces, and the
# two warning messages (blue "you might want to..." then red "couldn't find...")
# are reordered to both sit after the closing brace.
$p2 = $d.Paragraphs.Item(2)
$p2xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:r><w:t>{</w:t></w:r>' +
  '<w:r><w:t>m</w:t></w:r>' +
  '<w:r><w:t>:</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '<w:r><w:t>comment</w:t></w:r>' +
  '<w:r><w:t>block</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> some important comment</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">    </w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="0000FF"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>&lt;---</w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="0000FF"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>You might want to replace m: commentblock by m:commentblock</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">    </w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>&lt;---</w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>Couldn''t find the ''commentblock'' variable</w:t></w:r>' +
  '</w:p>'
$p2.Range.InsertXML($p2xml)

# Paragraph 4: the "m:endcommentblock" field paragraph.
# Same field -> literal-text-brace conversion; the trailing " m:endcommentblock "
# field-code run and the closing fldChar end run are dropped since the tag text
# is now folded into the single literal "{m:endcommentblock}" run.
$p4 = $d.Paragraphs.Item(4)
$p4xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:r><w:t xml:space="preserve">{m:endcommentblock}</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">    </w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>&lt;---</w:t></w:r>' +
  '<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:highlight w:val="lightGray"/></w:rPr><w:t>Unexpected tag m:endcommentblock at this location</w:t></w:r>' +
  '</w:p>'
$p4.Range.InsertXML($p4xml)
